$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 65

# Column A holds a date-looking string ("2025-05-02") that must stay a plain
# text value (as every other row in this sheet does), not get auto-converted
# to a numeric date serial by Excel's normal typed-input parsing. Temporarily
# force a text number format while the value is entered, then clear the
# formatting again so the cell ends up without any explicit style (matching
# the rest of the sheet).
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025-05-02"
$cellA.ClearFormats()

$ws.Cells.Item($row, 2).Value = "espèces exotiques envahissantes"
$ws.Cells.Item($row, 3).Value = 30
$ws.Cells.Item($row, 4).Value = 1
